# Append three new contact/proxy records to the "group_b" sheet.
# (commit: "added checkboxes for white list , space encoding, test email will be cached from now on")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - John Snow
$ws.Range("A5").Value = "John"
$ws.Range("B5").Value = "Snow"
$ws.Range("C5").Value = "john.snow42@outlook.hu"
$ws.Range("D5").Value = "Whateverpass0*"
$ws.Range("E5").Value = "81.28.96.148:4000"
$ws.Range("F5").Value = "i0BdGW79w6Oo"
$ws.Range("G5").Value = "5Ao37R1ry6bc"

# Row 6 - Louise Deforge
$ws.Range("A6").Value = "Louise"
$ws.Range("B6").Value = "Deforge"
$ws.Range("C6").Value = "saidunuhu579@gmail.com"
$ws.Range("D6").Value = "ujemlxfzbgwnkzpd"
$ws.Range("E6").Value = "185.125.171.221:4021"
$ws.Range("F6").Value = "2JjU2izT4rk1tGb"
$ws.Range("G6").Value = "bEjmSK36Ma4C36t"

# Row 7 - Janice Holley (reuses the same proxy as row 6)
$ws.Range("A7").Value = "Janice"
$ws.Range("B7").Value = "Holley"
$ws.Range("C7").Value = "lithbello60@gmail.com"
$ws.Range("D7").Value = "wrddcsqmfqdkqkrx"
$ws.Range("E7").Value = "185.125.171.221:4021"
$ws.Range("F7").Value = "2JjU2izT4rk1tGb"
$ws.Range("G7").Value = "bEjmSK36Ma4C36t"

# Leave the selection where the author left it when saving
$ws.Range("D17").Select()
